# Refresh the "need_to_buy" data table: shift the date series forward by one
# day and update the computed TISG/fcs/buy/MYDIR/need_to_buy figures for the
# visible 14-day window (rows 2-15).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45974
$ws.Cells.Item(2, 2).Value = 8843.33032346788
$ws.Cells.Item(2, 3).Value = 8245.69577563744
$ws.Cells.Item(2, 4).Value = 11328
$ws.Cells.Item(2, 5).Value = 4478.11739634064
$ws.Cells.Item(2, 6).Value = 58.1588821657531
$ws.Cells.Item(3, 1).Value = 45975
$ws.Cells.Item(3, 2).Value = 8695.17697736804
$ws.Cells.Item(3, 3).Value = 7301.09013198785
$ws.Cells.Item(3, 4).Value = 7200
$ws.Cells.Item(3, 5).Value = 4822.80809266108
$ws.Cells.Item(3, 6).Value = 205.162426027039
$ws.Cells.Item(4, 1).Value = 45976
$ws.Cells.Item(4, 2).Value = 2997.26622852637
$ws.Cells.Item(4, 3).Value = 4668.54309570319
$ws.Cells.Item(4, 4).Value = 7200
$ws.Cells.Item(4, 5).Value = 4287.33384944507
$ws.Cells.Item(4, 6).Value = 73.1615393811777
$ws.Cells.Item(5, 1).Value = 45977
$ws.Cells.Item(5, 2).Value = 2755.56763358653
$ws.Cells.Item(5, 3).Value = 4406.98070165649
$ws.Cells.Item(5, 4).Value = 7200
$ws.Cells.Item(5, 5).Value = 4152.75002116228
$ws.Cells.Item(5, 6).Value = 56.6554467841152
$ws.Cells.Item(6, 1).Value = 45978
$ws.Cells.Item(6, 2).Value = 9196.29931971738
$ws.Cells.Item(6, 3).Value = 7530.17138859087
$ws.Cells.Item(6, 4).Value = 7200
$ws.Cells.Item(6, 5).Value = 5157.44833100166
$ws.Cells.Item(6, 6).Value = 228.650821649688
$ws.Cells.Item(7, 1).Value = 45979
$ws.Cells.Item(7, 2).Value = 10116.508538002
$ws.Cells.Item(7, 3).Value = 8755.72435421698
$ws.Cells.Item(7, 4).Value = 7200
$ws.Cells.Item(7, 5).Value = 5829.54417612728
$ws.Cells.Item(7, 6).Value = 307.719522097678
$ws.Cells.Item(8, 1).Value = 45980
$ws.Cells.Item(8, 2).Value = 10116.508538002
$ws.Cells.Item(8, 3).Value = 9058.39321936113
$ws.Cells.Item(8, 4).Value = 7200
$ws.Cells.Item(8, 5).Value = 5829.54417612728
$ws.Cells.Item(8, 6).Value = 320.330724812017
$ws.Cells.Item(9, 1).Value = 45981
$ws.Cells.Item(9, 2).Value = 10116.508538002
$ws.Cells.Item(9, 3).Value = 9325.25777090924
$ws.Cells.Item(9, 4).Value = 7200
$ws.Cells.Item(9, 5).Value = 5829.54417612728
$ws.Cells.Item(9, 6).Value = 331.450081126522
$ws.Cells.Item(10, 1).Value = 45982
$ws.Cells.Item(10, 2).Value = 10116.508538002
$ws.Cells.Item(10, 3).Value = 8897.04134581315
$ws.Cells.Item(10, 4).Value = 7200
$ws.Cells.Item(10, 5).Value = 5829.54417612728
$ws.Cells.Item(10, 6).Value = 313.607730080851
$ws.Cells.Item(11, 1).Value = 45983
$ws.Cells.Item(11, 2).Value = 4000.30754528634
$ws.Cells.Item(11, 3).Value = 6923.97745169389
$ws.Cells.Item(11, 4).Value = 7200
$ws.Cells.Item(11, 5).Value = 5477.70511299437
$ws.Cells.Item(11, 6).Value = 216.736773528678
$ws.Cells.Item(12, 1).Value = 45984
$ws.Cells.Item(12, 2).Value = 3853.28814954241
$ws.Cells.Item(12, 3).Value = 7088.65804782756
$ws.Cells.Item(12, 4).Value = 7200
$ws.Cells.Item(12, 5).Value = 5469.70161357832
$ws.Cells.Item(12, 6).Value = 223.264985891912
$ws.Cells.Item(13, 1).Value = 45985
$ws.Cells.Item(13, 2).Value = 10980.8696790314
$ws.Cells.Item(13, 3).Value = 10548.7166241374
$ws.Cells.Item(13, 4).Value = 7200
$ws.Cells.Item(13, 5).Value = 6409.18197125632
$ws.Cells.Item(13, 6).Value = 406.579108141406
$ws.Cells.Item(14, 1).Value = 45986
$ws.Cells.Item(14, 2).Value = 10980.8696790314
$ws.Cells.Item(14, 3).Value = 10223.6172916093
$ws.Cells.Item(14, 4).Value = 7200
$ws.Cells.Item(14, 5).Value = 6409.18197125632
$ws.Cells.Item(14, 6).Value = 393.033302619401
$ws.Cells.Item(15, 1).Value = 45987
$ws.Cells.Item(15, 2).Value = 10980.8696790314
$ws.Cells.Item(15, 3).Value = 9881.82842779281
$ws.Cells.Item(15, 4).Value = 7200
$ws.Cells.Item(15, 5).Value = 6409.18197125632
$ws.Cells.Item(15, 6).Value = 378.79209996038
